$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8, 9, 10 (ECs-as-sending-cluster rows no longer present)
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fbln1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 80.365851
$ws.Range("H2").Value = 241.097553
$ws.Range("I2").Value = 0.9568768228420588
$ws.Range("J2").Value = 0.9568768228420588
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 6195.36514586769
$ws.Range("R2").Value = 55758.28631280921
$ws.Range("S2").Value = 0.2300156461450489
$ws.Range("T2").Value = 0.2300156461450489

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fbln1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 80.365851
$ws.Range("H3").Value = 241.097553
$ws.Range("I3").Value = 0.9568768228420588
$ws.Range("J3").Value = 0.9568768228420588
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 8163.566144905105
$ws.Range("R3").Value = 73472.09530414594
$ws.Range("S3").Value = 0.3030891476865818
$ws.Range("T3").Value = 0.3030891476865817

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fbln1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 80.365851
$ws.Range("H4").Value = 241.097553
$ws.Range("I4").Value = 0.9568768228420588
$ws.Range("J4").Value = 0.9568768228420588
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 11414.10378957106
$ws.Range("R4").Value = 102726.9341061395
$ws.Range("S4").Value = 0.4237720290104281
$ws.Range("T4").Value = 0.4237720290104281

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fbln1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.621815
$ws.Range("H5").Value = 10.865445
$ws.Range("I5").Value = 0.04312317715794126
$ws.Range("J5").Value = 0.04312317715794126
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 279.2039919515166
$ws.Range("R5").Value = 2512.83592756365
$ws.Range("S5").Value = 0.01036602122763349
$ws.Range("T5").Value = 0.01036602122763349

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fbln1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.621815
$ws.Range("H6").Value = 10.865445
$ws.Range("I6").Value = 0.04312317715794126
$ws.Range("J6").Value = 0.04312317715794126
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("Q6").Value = 367.9041029144266
$ws.Range("R6").Value = 3311.13692622984
$ws.Range("S6").Value = 0.01365919489147794
$ws.Range("T6").Value = 0.01365919489147794

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fbln1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.621815
$ws.Range("H7").Value = 10.865445
$ws.Range("I7").Value = 0.04312317715794126
$ws.Range("J7").Value = 0.04312317715794126
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("Q7").Value = 514.3947560093067
$ws.Range("R7").Value = 4629.55280408376
$ws.Range("S7").Value = 0.01909796103882983
$ws.Range("T7").Value = 0.01909796103882984
